$wb = $excel.ActiveWorkbook

# --- YDS sheet: append Week 17 play-by-play yardage logs ---
$ydsWs = $wb.Worksheets.Item("YDS")
$ydsWs.Range("B2").Value = '3 18 2 2 3 -3 0 4 1 16 3 2 4 11 6 -1 15 11 4 2 2 5 5 8 1 0 28 2 4 2 3 2 6 5 12 31 9 7 1 9 9 8 0 0 21 3 4 5 20 14 -4 3 4 -2 20 6 2 8 5 9 2 1 11 5 9 -1 3 1 1 2 2 -2 2 5 0 31 6 2 0 4 8 8 1 9 0 5 7 8 5 12 3 -1 7 5 2 1 1 3 4 5 7 11 2 1 -1 6 -2 3 4 3 5 4 0 2 3 5 3 5 5 4 5 0 1 3 3 2 -2 13 4 3 -1 9 2 4 0 12 5 -1 3 6 4 -3 6 6 6 2 8 6 16 6 1 11 14 3 3 22 2 3 5 5 7 3 6 1 5 4 1 7 2 2 0 -3 5 7 3 0 9 -1 9 3 2 12 2 10 -5 4 -2 -1 12 3 10 3 7 1 16 2 1 7 7 2 10 7 2 5 6 5 8 5 6 4 4 2 6 3 4 9 3 5 6 7 11 6 12 0 1 7 1 4 9 4 6 1 3 8 7 8 5 10 1 6 5 4 2 3 19 5 12 1 4 19 4 12 2 8 4 2 -2 0 0 5 2 11 0 -4 1 3 1 11 5 7 5 10 1 -1 -3 4 1 4 2 5 4 2 19 0 3 8 2 3 4 5 2 5 4 6 5 0 2 1 6 -1 2 2 3 1 2 4 1 1 6 1 0 3 4 0 13 9 5 4 2 0 5 1 0 4 2 4 1 4 3 2 8 13 1 3 11 8 3 1 4 6 5 4 1 12 5 18 -1 3 5 3 8 -2 4 5 3 8 6 8 0 -1 -2 13 3 2 2 3 1 5 1 4 5 4 3 -4 -1 7 2 18 4 8 5 0 11 13 21 5 3 2 1 1 2 4 1 7 2 -2 5 3 16 4 6 8 12 6 13 4 -1 6 10 7 3 15 2 8 8 2 6 1 0 7 -2 3 -5 3 4 2 2 -7 1 2 1 2 -3 5 3 7 1 1 1 9 5 17 11 3 -1 12 2 4 3 3 32 6 6 4 6 3 2 3 -2'
$ydsWs.Range("B3").Value = '2 8 0 21 6 6 29 29 10 13 4 10 6 4 13 6 49 9 10 11 17 -4 7 20 8 9 7 11 4 4 12 8 27 20 42 19 9 10 6 5 16 11 6 21 37 19 29 3 9 41 7 24 5 18 20 6 12 24 12 9 49 32 7 13 14 4 8 13 7 8 17 4 7 14 32 4 5 5 17 15 28 7 28 2 6 11 12 8 8 7 17 3 19 7 43 27 18 7 6 15 5 17 10 9 20 8 4 14 4 11 7 7 5 7 4 4 8 16 8 9 21 1 10 6 4 11 9 8 20 13 6 1 24 6 23 19 20 18 7 1 35 39 25 3 13 13 -1 10 12 7 11 13 13 15 7 -1 9 4 5 12 22 3 3 8 8 18 11 5 22 6 8 7 7 20 20 7 13 7 4 11 15 7 5 6 13 3 11 9 9 11 0 8 7 8 4 8 20 5 8 30 6 9 15 11 1 16 5 9 8 -2 11 11 2 3 12 9 15 2 9 8 -2 3 9 4 12 7 21 6 29 8 2 7 5 6 9 8 13 8 8 4 7 39 13 11 5 -1 0 10 3 14 16 18 0 9 9 8 2 10 29 3 10 15 15 3 17 6 5 22 18 13 5 6 11 2 -1 5 3 5 4 9 4 20 10 32 4 8 -3 20 36 12 3 11 2 6 5 6 4 8 14 3 43 5 5 4 22 8 4 8 7 9 2 7 1 9 12 5 5 3 5 5 5 3 3 6 12 3 11 17 5 15 6 4 9 -1 16 12 10 -1 18 28 11 11 4 7 12 17 25 18 14 4 6 7 14 5 15 13 5 1 14 9 18 5 33 2 14 15 4 15 7 11 4 5'
$ydsWs.Range("C2").Value = '31 -1 2 8 0 13 0 0 0 2 2 2 4 1 15 3 -2 6 0 1 9 1 7 -1 1 3 9 4 2 7 6 15 1 3 -5 4 -2 -2 11 1 4 14 2 6 4 14 4 0 9 -2 2 1 2 3 12 4 0 4 0 1 1 -1 0 1 -4 7 8 1 4 14 3 2 31 3 3 3 3 -3 6 10 -1 12 3 11 4 4 3 13 2 2 -1 2 2 19 13 0 4 4 8 22 -1 7 2 1 -1 3 -4 1 -1 4 2 5 2 10 -2 2 1 3 -2 2 9 0 -1 -1 9 0 0 5 1 0 0 10 7 21 3 46 3 -3 -5 2 2 1 11 15 -3 4 66 1 3 3 0 0 1 2 -1 1 0 -1 0 -8 1 24 1 3 -1 3 1 -1 5 4 2 4 2 1 -1 4 1 14 0 11 0 1 7 2 1 8 9 2 11 5 -6 3 10 7 5 6 4 5 5 -1 6 11 2 0 -3 4 1 2 7 2 2 -1 3 4 -1 5 1 9 3 -3 5 2 7 0 5 -1 5 4 2 3 6 2 2 0 1 8 3 2 0 8 13 3 5 8 -2 2 1 3 3 2 5 1 6 15 1 2 13 -2 5 6 3 4 11 8 11 1 2 7 1 0 -2 1 4 2 2 4 2 3 7 9 8 -1 4 11 5 4 9 3 5 5 6 0 1 0 5 8 -2 19 3 1 1 9 1 1 1 3 12 2 2 2 8 2 2 2 6 2 4 1 1 2 3 14 1 7 7 -2 11 4 7 1 8 1 -4'
$ydsWs.Range("C3").Value = '9 13 6 5 6 12 9 24 8 15 13 7 4 12 16 2 3 4 21 4 5 16 9 7 9 37 10 10 20 18 11 27 32 31 6 7 19 10 33 17 11 14 10 4 14 20 14 5 11 40 5 46 3 8 3 23 13 7 11 3 4 5 8 1 15 10 5 12 11 14 19 9 2 -3 11 9 6 24 19 22 11 19 3 3 12 14 3 11 11 8 3 5 6 13 -3 9 9 32 9 76 11 15 6 8 15 11 7 -2 5 9 16 -1 31 42 21 24 18 4 4 29 23 2 15 11 11 8 3 18 7 5 8 26 1 16 1 5 10 21 15 -1 5 4 5 10 6 11 23 6 5 9 55 7 4 6 13 26 7 17 21 32 7 82 25 2 11 11 27 6 50 5 3 7 7 20 4 18 5 9 5 6 27 1 2 16 9 11 7 21 2 6 19 2 52 11 9 1 6 14 35 23 6 64 6 29 22 22 23 60 10 2 3 12 9 4 23 11 9 36 10 0 -1 41 11 7 4 16 6 38 8 20 16 8 3 18 7 12 4 17 5 14 7 7 4 9 40 29 25 9 2 5 9 5 5 -1 8 3 9 4 34 5 5 9 7 6 27 9 1 12 5 7 13 17 3 3 4 17 6 11 31 7 22 10 3 13 3 24 6 1 16 9 7 6 9 25 14 11 23 -6 10 18 20 3 11 16 5 17 3 68 20 17 9 11 15 3 52 1 19 3 18 9 3 2 5 14 13 0 12 4 22 3 10 20 4 52 10 18 7 19 6 4 9 22 11 17 11 18 2 35 12 9 2 8 3 21 15 24 4 15 7'

# --- OFF sheet: Home (row 2) and Road (row 3) season totals ---
$offWs = $wb.Worksheets.Item("OFF")
$offWs.Range("C2").Value = 205
$offWs.Range("D2").Value = 14
$offWs.Range("F2").Value = 60
$offWs.Range("G2").Value = 74
$offWs.Range("J2").Value = 45
$offWs.Range("L2").Value = 282
$offWs.Range("M2").Value = 190
$offWs.Range("O2").Value = 25
$offWs.Range("P2").Value = 17
$offWs.Range("Q2").Value = 581
$offWs.Range("B3").Value = 10
$offWs.Range("C3").Value = 217
$offWs.Range("E3").Value = 30
$offWs.Range("F3").Value = 126
$offWs.Range("G3").Value = 31
$offWs.Range("H3").Value = 32
$offWs.Range("I3").Value = 62
$offWs.Range("J3").Value = 40
$offWs.Range("N3").Value = 24

# --- DEF sheet: Home (row 2) and Road (row 3) season totals ---
$defWs = $wb.Worksheets.Item("DEF")
$defWs.Range("C2").Value = 173
$defWs.Range("F2").Value = 40
$defWs.Range("G2").Value = 50
$defWs.Range("J2").Value = 22
$defWs.Range("L2").Value = 274
$defWs.Range("M2").Value = 178
$defWs.Range("O2").Value = 17
$defWs.Range("P2").Value = 10
$defWs.Range("Q2").Value = 453
$defWs.Range("B3").Value = 16
$defWs.Range("C3").Value = 185
$defWs.Range("E3").Value = 37
$defWs.Range("F3").Value = 102
$defWs.Range("G3").Value = 34
$defWs.Range("H3").Value = 31
$defWs.Range("I3").Value = 63
$defWs.Range("J3").Value = 46
$defWs.Range("N3").Value = 12

# --- ST sheet: kicking/special-teams totals + per-week logs ---
$stWs = $wb.Worksheets.Item("ST")
$stWs.Range("B2").Value = 83
$stWs.Range("D2").Value = 62
$stWs.Range("F2").Value = 385
$stWs.Range("G2").Value = 381
$stWs.Range("J2").Value = 180
$stWs.Range("K2").Value = 177
$stWs.Range("L2").Value = 110
$stWs.Range("M2").Value = 99
$stWs.Range("B3").Value = 46
$stWs.Range("B4").Value = '63 65 57 65 55 48 67 58 63 64 65 61 63 60 67 59 67 63 61 39 65 61 64 58 63 63 47 40 65 63 65 53 56 59 59 43 59'
$stWs.Range("B5").Value = '13 20 8 33 3 9 25 17 19 26 20 16 18 16 36 10 23 98 17 0 15 14 19 13 21 0 13 6 35 21 20 0 9 25 20 0 18'
$stWs.Range("B6").Value = '23 23 21 23 27 23 47 21 24 24 21 25 25 26 35 0 19 12 23 24 26 34 21 22 15 24 0 22'
$stWs.Range("D3").Value = '47 42 49 41 49 46 54 48 35 34 52 55 46 49 57 42 49 53 47 49 53 33 48 49 40 31 50 46 39 35 51 45 37 53 43 39 24 41 39 37 42 44 35 57 45 53 33 50 55 45 46 57 60 44 44 42 47 43 43 36 42 61'
$stWs.Range("D4").Value = '3 0 0 7 0 13 11 0 0 0 0 5 8 0 9 0 6 0 13 11 0 0 13 7 0 0 0 0 0 0 0 6 0 0 0 0 0 4 0 0 0 0 0 4 0 0 0 -3 14 5 6 8 8 1 6 0 6 10 0 0 0 18'
$stWs.Range("D5").Value = '0 22 15 30 0 0 16 0 0 0 0 13 29 0 12 0 0 10 0 42 0 0 0 0 0 7 0 14 13 17 0 0 16 21 0 0 0 0 14 0 0 0 0 0 9 0 0 0 7 12 0 0 0 13 22 0 0 1 0 0 0 12 0 0 0 0 0 0 0 0 0 0 0 0'

# --- TURNS sheet: Home (row 2) and Road (row 3) turnover totals ---
$turnsWs = $wb.Worksheets.Item("TURNS")
$turnsWs.Range("B2").Value = 12
$turnsWs.Range("C2").Value = 5
$turnsWs.Range("E2").Value = 5
$turnsWs.Range("D3").Value = 5

# --- PEN sheet: penalty counts ---
$penWs = $wb.Worksheets.Item("PEN")
$penWs.Range("B2").Value = 17
$penWs.Range("B3").Value = 27

